$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 378.925108
$ws.Range("D2").Value = 1820.247081

$ws.Range("B3").Value = 6.533915
$ws.Range("D3").Value = 15.693523
$ws.Range("E3").Value = 0

$ws.Range("B4").Value = 69.113219
$ws.Range("C4").Value = 332

$ws.Range("G5").Value = 0.009835999999999999
$ws.Range("H5").Value = -0.155587
$ws.Range("I5").Value = 0.175259
$ws.Range("J5").Value = 0.989255

$ws.Range("G6").Value = -0.286166
$ws.Range("H6").Value = -0.459916
$ws.Range("I6").Value = -0.112416
$ws.Range("J6").Value = 0.000373

$ws.Range("G7").Value = -0.296002
$ws.Range("H7").Value = -0.426651
$ws.Range("I7").Value = -0.165354
$ws.Range("J7").Value = 0.000001
